# "Add today's walk around Athens"
# The workbook tracks a running daily total in column F ("Sheet1"),
# computed from that day's distance entered in column G
# (F2 = F1 + G2). Today's walk adds 2.3 to the running total, i.e.
# G2 (today's distance) goes from 56.9 to 59.2, which in turn
# recalculates F2 (running total) from 106.9 to 109.2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("G2").Value = 59.2
